# Daily "tick" update: advance the schedule by one day.
# For every data row (2..99):
#   - E (剩余 / remaining days) decreases by 1.
#   - If that would bring E down to 0 (i.e. the cycle finished), the row is
#     reset: E goes back to D (总天 / total days) and F (开始时间 / start date)
#     is bumped to the new "today" (2025-11-18).
#   - Otherwise F is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newToday = 20251118

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # column D
    $eCell = $ws.Cells.Item($r, 5)   # column E
    $fCell = $ws.Cells.Item($r, 6)   # column F

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null) {
        continue
    }

    # Skip rows whose start-date (F, YYYYMMDD) isn't a well-formed 8-digit
    # date - e.g. a corrupted value like 202510929 - mirroring the source
    # update tool, which leaves such rows untouched rather than failing.
    $fText = [string]$fVal
    if ($fText.Length -ne 8) {
        continue
    }

    $newE = $eVal - 1

    if ($newE -le 0) {
        $eCell.Value2 = $dVal
        $fCell.Value2 = $newToday
    } else {
        $eCell.Value2 = $newE
    }
}
